# Update the "Generate Report for Handback" timestamps/status across sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for rows 2 and 5
$wsOverview.Range("G2").Value = "2016-09-04 16:19:00"
$wsOverview.Range("G5").Value = "2016-09-04 16:19:00"

# zh-cn sheet: Priority (E) changes from "ht" to "mt"; Correspond Handoff
# Datetime (H) and Correspond Handback DateTime (K) get refreshed timestamps.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-04 16:18:56"
$wsZhCn.Range("H5").Value = "2016-09-04 16:18:56"
$wsZhCn.Range("K2").Value = "2016-09-04 16:19:16"
$wsZhCn.Range("K5").Value = "2016-09-04 16:19:16"

# de-de sheet: Priority (E) changes from "ht" to "mt"; Correspond Handoff
# Datetime (H) and Correspond Handback DateTime (K) get refreshed timestamps.
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-09-04 16:19:00"
$wsDeDe.Range("H5").Value = "2016-09-04 16:19:00"
$wsDeDe.Range("K2").Value = "2016-09-04 16:19:24"
$wsDeDe.Range("K5").Value = "2016-09-04 16:19:24"
